# Work time log ("Työaikakirjanpito") - add a new entry on row 58
# and reselect cell C46 (matches the author's recorded end state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry: date 2020-09-08 (serial 44082), 1 hour,
# "Frontin jakamista ja backendin malleihin lisäystä"
$ws.Range("A58").Value = 44082
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = "Frontin jakamista ja backendin malleihin lisäystä"

# The SUM formula in B62 (=SUM(B$2:B$61)) recalculates automatically to 184.5.

# Restore the selection the author left the sheet in.
$ws.Range("C46").Select()
